# Zeiterfassung correction (Milestone 2):
#   - "Geplant"     (row 27, col C): 40 Stunden -> 50 Stunden
#   - "Mehraufwand" (row 28, col C): 2 Stunden 15 Minuten -> Keiner
#   - switch calculation back to automatic (was left in manual mode)
#   - leave the selection on the last touched cell (C29)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

$ws.Range("C27").Value = "50 Stunden"
$ws.Range("C28").Value = "Keiner"

$excel.Calculation = -4105   # xlCalculationAutomatic

$ws.Range("C29").Select() | Out-Null
